$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 31 (R5 / All Industries row)
$ws.Rows.Item(31).Delete()

# Add column D header and values
$dValues = @(1,1,2,3,4,5,5,5,6,5,5,7,7,7,7,8,8,8,8,8,8,8,9,9,9,9,3,10,6,10)
for ($i = 0; $i -lt $dValues.Length; $i++) {
    $row = $i + 1
    $cell = $ws.Cells.Item($row, 4)
    $cell.Value = $dValues[$i]
    $cell.NumberFormat = "General"
}
